$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new blank columns before column D ---
# (shifts the old D:K quarterly columns to F:M, carrying formats along)
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Populate the two new columns with the newly reported quarter figures ---
$ws.Range("D7").Value = 43465; $ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 119500; $ws.Range("E8").Value = 120300
$ws.Range("D9").Value = 25300; $ws.Range("E9").Value = 24600
$ws.Range("D10").Value = 94200; $ws.Range("E10").Value = 95700
$ws.Range("D12").Value = 6200; $ws.Range("E12").Value = 5400
$ws.Range("D13").Value = 0; $ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0; $ws.Range("E14").Value = 0
$ws.Range("D15").Value = 33700; $ws.Range("E15").Value = 33400
$ws.Range("D17").Value = 102300; $ws.Range("E17").Value = 98800
$ws.Range("D18").Value = 17200; $ws.Range("E18").Value = 21500
$ws.Range("D20").Value = -20400; $ws.Range("E20").Value = -20100
$ws.Range("D21").Value = 30500; $ws.Range("E21").Value = 34700
$ws.Range("D22").Value = 0; $ws.Range("E22").Value = 0
$ws.Range("D23").Value = -3200; $ws.Range("E23").Value = 1400
$ws.Range("D24").Value = -7200; $ws.Range("E24").Value = -700
$ws.Range("D25").Value = 0; $ws.Range("E25").Value = 0
$ws.Range("D26").Value = 4000; $ws.Range("E26").Value = 2100
$ws.Range("D27").Value = 4000; $ws.Range("E27").Value = 2100
$ws.Range("D28").Value = 0; $ws.Range("E28").Value = 0
$ws.Range("D29").Value = 200; $ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0; $ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0; $ws.Range("E31").Value = 0
$ws.Range("D32").Value = 20400; $ws.Range("E32").Value = 20100
$ws.Range("D33").Value = 4200; $ws.Range("E33").Value = 2100
$ws.Range("D34").Value = 0; $ws.Range("E34").Value = 0
$ws.Range("D35").Value = 4200; $ws.Range("E35").Value = 2100
$ws.Range("D38").Value = 43465; $ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 297500; $ws.Range("E41").Value = 128700
$ws.Range("D42").Value = 0; $ws.Range("E42").Value = 0
$ws.Range("D43").Value = 146700; $ws.Range("E43").Value = 286900
$ws.Range("D44").Value = 24400; $ws.Range("E44").Value = 25900
$ws.Range("D45").Value = 20300; $ws.Range("E45").Value = 22700
$ws.Range("D46").Value = 488900; $ws.Range("E46").Value = 464200
$ws.Range("D47").Value = 8800; $ws.Range("E47").Value = 8800
$ws.Range("D48").Value = 116300; $ws.Range("E48").Value = 118000
$ws.Range("D49").Value = 927900; $ws.Range("E49").Value = 936900
$ws.Range("D50").Value = 0; $ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0; $ws.Range("E51").Value = 0
$ws.Range("D52").Value = 6300; $ws.Range("E52").Value = 6300
$ws.Range("D53").Value = 0; $ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1548300; $ws.Range("E54").Value = 1534200
$ws.Range("D57").Value = 421700; $ws.Range("E57").Value = 393000
$ws.Range("D58").Value = 8200; $ws.Range("E58").Value = 8200
$ws.Range("D59").Value = 41700; $ws.Range("E59").Value = 51500
$ws.Range("D60").Value = 471600; $ws.Range("E60").Value = 452700
$ws.Range("D61").Value = 1155000; $ws.Range("E61").Value = 1156200
$ws.Range("D62").Value = 30500; $ws.Range("E62").Value = 38500
$ws.Range("D63").Value = 0; $ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0; $ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0; $ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1657200; $ws.Range("E66").Value = 1647500
$ws.Range("D68").Value = 0; $ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0; $ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0; $ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0; $ws.Range("E71").Value = 0
$ws.Range("D72").Value = -229500; $ws.Range("E72").Value = -233700
$ws.Range("D73").Value = 0; $ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0; $ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0; $ws.Range("E75").Value = 0
$ws.Range("D76").Value = -108900; $ws.Range("E76").Value = -113200
$ws.Range("D77").Value = 0; $ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465; $ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 4200; $ws.Range("E81").Value = 2100
$ws.Range("D83").Value = 33700; $ws.Range("E83").Value = 33400
$ws.Range("D84").Value = 0; $ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0; $ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0; $ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0; $ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0; $ws.Range("E88").Value = 0
$ws.Range("D89").Value = 201900; $ws.Range("E89").Value = 43200
$ws.Range("D91").Value = -24500; $ws.Range("E91").Value = -20600
$ws.Range("D92").Value = 0; $ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0; $ws.Range("E93").Value = 0
$ws.Range("D94").Value = -29600; $ws.Range("E94").Value = -25800
$ws.Range("D96").Value = 0; $ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0; $ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0; $ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0; $ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2000; $ws.Range("E100").Value = 1100
$ws.Range("D101").Value = -900; $ws.Range("E101").Value = 200
$ws.Range("D102").Value = 169300; $ws.Range("E102").Value = 18600

# --- Minor restatements to a few previously reported quarters ---
$ws.Range("H89").Value = 26300
$ws.Range("I89").Value = -4100
$ws.Range("I91").Value = -26400
$ws.Range("J91").Value = -26500
$ws.Range("I94").Value = -36400
$ws.Range("H102").Value = 20400
$ws.Range("I102").Value = -40000
